$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New timestamp rows to append below the existing data (rows 2-39 already populated).
$timestamps = @(
    "20 Jan 2024 - 18:44 UTC",
    "20 Jan 2024 - 18:44 UTC",
    "20 Jan 2024 - 18:47 UTC",
    "20 Jan 2024 - 18:48 UTC",
    "20 Jan 2024 - 18:49 UTC",
    "20 Jan 2024 - 18:50 UTC",
    "20 Jan 2024 - 18:51 UTC",
    "20 Jan 2024 - 18:52 UTC",
    "20 Jan 2024 - 18:53 UTC",
    "20 Jan 2024 - 18:54 UTC",
    "20 Jan 2024 - 18:55 UTC",
    "20 Jan 2024 - 18:56 UTC",
    "20 Jan 2024 - 18:57 UTC",
    "20 Jan 2024 - 18:58 UTC",
    "20 Jan 2024 - 18:59 UTC",
    "20 Jan 2024 - 19:00 UTC",
    "20 Jan 2024 - 19:01 UTC",
    "20 Jan 2024 - 19:02 UTC",
    "20 Jan 2024 - 19:03 UTC",
    "20 Jan 2024 - 19:04 UTC",
    "20 Jan 2024 - 19:05 UTC",
    "20 Jan 2024 - 19:06 UTC"
)

$startRow = 40
for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $timestamps[$i]
    $ws.Cells.Item($row, 2).Value = 0
}
